$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the Job_Description text for the new row by repeating the
# existing JD_001 description text 5 times (matches source data exactly).
$desc1 = "We are seeking a Software Engineer to build and maintain high-quality software solutions." + "`n" + "Work with global teams to drive innovation and deliver scalable applications." + "`n" + "Join Akkodis and be part of a tech-driven, collaborative environment."
$newDesc = $desc1 + $desc1 + $desc1 + $desc1 + $desc1

# Add new row 7 with JD_006 posting data
$ws.Cells.Item(7, 1).Value = "JD_006"
$ws.Cells.Item(7, 2).Value = "Senior X Engineer"
$ws.Cells.Item(7, 3).Value = $newDesc
$ws.Cells.Item(7, 4).Value = 2
$ws.Cells.Item(7, 5).Value = 4

# The embedded newlines in the description cause the engine to auto
# expand the row height; AutoFit it back so the row keeps using the
# sheet's default height (matching the other data rows).
$ws.Rows.Item(7).EntireRow.AutoFit()
